# Fix: decimal point had been replaced by a comma when these figures were
# first entered, so every quarterly figure ended up stored as text in
# sharedStrings (e.g. "6.83", "19.3", ...). Re-enter the columns C:F for
# rows 2-6 as real numbers (fixing the truncated/garbled values along the
# way), flag the one genuinely unavailable figure as "NA", and make sure
# the corrected numbers are still right-aligned with an explicit black font
# like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- corrected numeric data (rows 2-6, columns C-F) -------------------
$ws.Range("C2").Value = 6.83
$ws.Range("D2").Value = 6.8
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = 6.71

$ws.Range("C3").Value = 23.34
$ws.Range("D3").Value = 23.34
$ws.Range("E3").Value = 23.2
$ws.Range("F3").Value = 23.68

$ws.Range("C4").Value = 19.35
$ws.Range("D4").Value = 19.38
$ws.Range("E4").Value = 19.58
$ws.Range("F4").Value = 19.68

$ws.Range("C5").Value = 3.84
$ws.Range("D5").Value = 3.84
$ws.Range("E5").Value = 3.54
$ws.Range("F5").Value = 3.74

$ws.Range("C6").Value = 19.38
$ws.Range("D6").Value = 19.58
$ws.Range("E6").Value = 19.35
$ws.Range("F6").Value = 19.69

# --- formatting: keep the values right aligned with an explicit font color
$data = $ws.Range("C2:F6")
$data.HorizontalAlignment = -4152
$data.Font.Color = 0
$data.Font.Name = "Calibri"

# --- selection left the way the editor left it after fixing the block --
$ws.Range("C2:F6").Select() | Out-Null
